$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell D holds numeric-looking text (e.g. "27.183.87", "1.001") that must
# stay plain text, exactly as authored upstream (t="inlineStr"/t="s").
# Setting .Value directly on a General-formatted cell would make Excel
# auto-coerce it to a number (dropping trailing zeros / misreading the
# multi-dot values), so we briefly force Text format, assign the literal
# string, then restore the Normal style so no stray style index is left
# behind (keeps cell formatting identical to the original).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.183.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5249"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3773"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07272"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8984"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08395"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.905.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.270"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.222.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.060"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.133.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.437"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.278"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.758"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.933"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.794"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09289"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05061"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.237"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.958"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.348"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.599"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5711"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01983"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.071"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.670"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1515"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4837"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.617"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.38%  "
